# Auto-generated edit script for LOQ4037.docx reorder
$d = $word.ActiveDocument
$brk = [char]11

# --- Paragraph 6 ---
$pStart = $d.Paragraphs.Item(6).Range.Start
$pEnd = $d.Paragraphs.Item(6).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = 'Propriedade gerais dos compostos orgânicos. Estrutura, métodos de obtenção, propriedades físicas, reações dos hidrocarbonetos alifáticos e aromáticos, haletos orgânicos, álcoois e características estruturais'+$brk+'como estereoquímica e a relação estrutura-reatividade.'

# --- Paragraph 7 ---
$pStart = $d.Paragraphs.Item(7).Range.Start
$pEnd = $d.Paragraphs.Item(7).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = 'General property of organic compounds. Physical properties, reactions of aliphatic and aromatic hydrocarbons, organic halides, ethers, alcohols and structural characteristics as stereochemistry and structure-reactivity.'
$fStart = $d.Paragraphs.Item(7).Range.Start
$fEnd = $d.Paragraphs.Item(7).Range.End
$fr = $d.Range($fStart, $fEnd)
$fr.Font.Italic = 1

# --- Paragraph 9 ---
$pStart = $d.Paragraphs.Item(9).Range.Start
$pEnd = $d.Paragraphs.Item(9).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = 'Gerais - Apresentar e Ensinar conceitos de Química Orgânica como instrumentos importantes para a compreensão de estratégias e operações industriais e tecnológicas. Abordar problemáticas sociais e ambientais com as quais a engenharia química está relacionada, tornando-os dessa forma, aptos a exercerem a função de Engenheiro Químico, e realizarem as mudanças que se façam necessárias.'+$brk+''+$brk+'Específicos – Compreender e descrever o mecanismo das reações orgânicas e a sua importância para o aprimoramento e desenvolvimento de processos industriais sintéticos e de etapas de formulação. Aprofundar o conceito de estrutura-reatividade e propriedades dos materiais.'

# --- Paragraph 11 ---
$pStart = $d.Paragraphs.Item(11).Range.Start
$pEnd = $d.Paragraphs.Item(11).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = '1.Teoria de Bronsted e de Lewis e acidez de compostos orgânicos'+$brk+'2.Alcanos - Processos de obtenção, Propriedades físicas, Análise Conformacional. Reação de Substituição Radicalar. '+$brk+'3.Isomeria Constitucional e Isomeria Espacial (Estereoquímica). Quiralidade, Nomenclatura R/S, classificação de estereoisômeros. Polarímetro e Técnicas de  Resolução de Isômeros Espaciais.'+$brk+'4.Haletos de Alquila – Substituição Nucleofílica, SN1, SN2, E1, E2. '+$brk+'5.Alcenos, Alcadienos e Alcinos – Propriedades físicas e químicas. Reação de adição eletrofílica (hidroalogenação, Hidratação, Halogenação, Diels-Alder, Redução-Oxidação). Adição conjugada em dienos (produto termodinâmico e cinético) '+$brk+'6. Fundamentos de RMN, Infra-vermelho, Ultra-violeta e Fluorescencia '+$brk+'7.Compostos aromáticos – Propriedades físicas dos aromáticos. Reações de Substituição Eletrofílica Aromática. Efeito de Grupos Substituintes. Reação de Substituição Nucleofílica.'+$brk+'8.Álcoois e Éteres – Propriedades físicas, reações e mecanismos.'

# --- Paragraph 12 ---
$pStart = $d.Paragraphs.Item(12).Range.Start
$pEnd = $d.Paragraphs.Item(12).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = 'Overview - Introduce and teach concepts of organic chemistry as important tools for understanding strategies and industrial and technological operations. Address social and environmental issues with which chemical engineering is related, making them thus able to exercise Chemical Engineer function, and realize the changes that are necessary.'+$brk+''+$brk+'Specific - Understand and describe the mechanism of organic reactions and their importance to the improvement and development of synthetic manufacturing processes and formulation stages. Deepening the concept of structure-reactivity and properties of materials.'
$fStart = $d.Paragraphs.Item(12).Range.Start
$fEnd = $d.Paragraphs.Item(12).Range.End
$fr = $d.Range($fStart, $fEnd)
$fr.Font.Italic = 1

# --- Paragraph 14 ---
$pStart = $d.Paragraphs.Item(14).Range.Start
$pEnd = $d.Paragraphs.Item(14).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = 'Duas provas teóricas e ao longo do semestre letivo'+$brk+''+$brk+'Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'

# --- Paragraph 19 ---
$pStart = $d.Paragraphs.Item(19).Range.Start
$pEnd = $d.Paragraphs.Item(19).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = '210064 - Eduardo Rezende Triboni'

# --- Paragraph 17 (multi-run: bold labels + plain content) ---
$pStart = $d.Paragraphs.Item(17).Range.Start
$pEnd = $d.Paragraphs.Item(17).Range.End
$r = $d.Range($pStart, $pEnd)
$r.Text = 'Método: A média final (M) será calculada pela expressão: M = (P1 + P2)/2'+$brk+'Critério: Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada recuperação com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'+$brk+'Norma de recuperação: BRESLOW, R. Questões e Exercícios de Química Orgânica. São Paulo: Makrons Books Editora, 1996. '+$brk+''+$brk+'BRUICE, P. Y. Química Orgânica, vol 1 e 2, São Paulo: Editora Pearson Prentice Hall, 2006. '+$brk+' '+$brk+'HENDRIKSON, James B.; CRAM, Donald J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966.'+$brk+''+$brk+'MCMURRY, John. Química Orgânica. São Paulo: Editora Pioneira Thomson Leraning, 2005.'+$brk+''+$brk+'SOLOMONS, T.W.G; FRYHLE, Graig. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001.'+$brk+''+$brk+'MORRISON, R.; BOYD, R. Química Orgânica. São Paulo: Editora Calouste Gulbenkian, 2008.'
$segStart = $d.Paragraphs.Item(17).Range.Start + 0
$segEnd = $d.Paragraphs.Item(17).Range.Start + 8
$sr = $d.Range($segStart, $segEnd)
$sr.Font.Bold = 1
$segStart = $d.Paragraphs.Item(17).Range.Start + 8
$segEnd = $d.Paragraphs.Item(17).Range.Start + 73
$segStart = $d.Paragraphs.Item(17).Range.Start + 73
$segEnd = $d.Paragraphs.Item(17).Range.Start + 83
$sr = $d.Range($segStart, $segEnd)
$sr.Font.Bold = 1
$segStart = $d.Paragraphs.Item(17).Range.Start + 83
$segEnd = $d.Paragraphs.Item(17).Range.Start + 355
$segStart = $d.Paragraphs.Item(17).Range.Start + 355
$segEnd = $d.Paragraphs.Item(17).Range.Start + 377
$sr = $d.Range($segStart, $segEnd)
$sr.Font.Bold = 1
$segStart = $d.Paragraphs.Item(17).Range.Start + 377
$segEnd = $d.Paragraphs.Item(17).Range.Start + 959

Write-Output "Paragraphs count: $($d.Paragraphs.Count)"